$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Translate the benchmark labels (column A) to Spanish ---
$ws.Range("A1").Value = "Marcos de Prueba (Intentos)"
$ws.Range("A2").Value = "MMLU (5-intentos)"
$ws.Range("A3").Value = "TriviaQA (1-intento)"
$ws.Range("A4").Value = "Preguntas Naturales (1-intento)"
$ws.Range("A5").Value = "QSM8K (8-intentos)"
$ws.Range("A6").Value = "EvalHumana (0-intentos)"
$ws.Range("A7").Value = "BBH (3-intentos)"

# --- 2. Add a thin box border around the whole table (A1:F7) ---
$table = $ws.Range("A1:F7")
$table.Borders.LineStyle = 1
$table.Borders.Weight = 2

# Keep the existing centered alignment explicit on the table range
$table.HorizontalAlignment = -4108
$table.VerticalAlignment = -4108

# --- 3. Wrap the long "Preguntas Naturales (1-intento)" label and grow its row ---
$ws.Range("A4").WrapText = $true
$ws.Rows(4).RowHeight = 28.5

# --- 4. Widen column A so the longer Spanish labels fit better ---
$ws.Range("A1").EntireColumn.ColumnWidth = 24.333333333333332
